# Update the Ptn -> Ptprz1 ligand-receptor table with the new TPM-derived
# NATMI statistics. The sending/target cluster grid is expanded from a
# partial 4x3 matrix to the full 4x4 matrix (ECs, FAPs, MuSCs,
# Resolving-Mac x ECs, FAPs, MuSCs, Resolving-Mac).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a numeric value while working around the script engine's
# lack of support for bare scientific-notation literals (e.g. 1.2E-05).
function Set-Num($cellRef, $text) {
    $ws.Range($cellRef).Value2 = [double]$text
}

# Each entry: SendingCluster, TargetCluster, then E..T (16 numeric values)
$rows = @(
    @("ECs","ECs","3","1","0.7169949999999999","2.150985","0.003651663653539308","0.003651663653539308","2","0.6666666666666666","0.138766","0.416298","0.01356925767068476","0.01356925767068476","0.09949452816999999","0.8954507535299999","4.955036504154898E-05","4.955036504154899E-05"),
    @("ECs","FAPs","3","1","0.7169949999999999","2.150985","0.003651663653539308","0.003651663653539308","2","0.6666666666666666","0.05416133333333333","0.162484","0.00529617548814441","0.005296175488144411","0.03883340519333333","0.34950064674","1.933985153282275E-05","1.933985153282275E-05"),
    @("ECs","MuSCs","3","1","0.7169949999999999","2.150985","0.003651663653539308","0.003651663653539308","3","1","10.002366","30.007098","0.9780831152479456","0.9780831152479456","7.171646410169999","64.54481769153","0.003571630562091421","0.003571630562091421"),
    @("ECs","Resolving-Mac","3","1","0.7169949999999999","2.150985","0.003651663653539308","0.003651663653539308","1","0.3333333333333333","0.03120566666666667","0.09361700000000001","0.003051451593225274","0.003051451593225274","0.02237430697166667","0.201368762745","1.114287487351535E-05","1.114287487351535E-05"),

    @("FAPs","ECs","3","1","191.4773253333333","574.431976","0.9751961860217362","0.9751961860217361","2","0.6666666666666666","0.138766","0.416298","0.01356925767068476","0.01356925767068476","26.57054252720533","239.134882744848","0.01323268832759796","0.01323268832759796"),
    @("FAPs","FAPs","3","1","191.4773253333333","574.431976","0.9751961860217362","0.9751961860217361","2","0.6666666666666666","0.05416133333333333","0.162484","0.00529617548814441","0.005296175488144411","10.37066724315378","93.33600518838399","0.005164810136540236","0.005164810136540236"),
    @("FAPs","MuSCs","3","1","191.4773253333333","574.431976","0.9751961860217362","0.9751961860217361","3","1","10.002366","30.007098","0.9780831152479456","0.9780831152479456","1915.226288685072","17237.03659816565","0.9538229236020548","0.9538229236020547"),
    @("FAPs","Resolving-Mac","3","1","191.4773253333333","574.431976","0.9751961860217362","0.9751961860217361","1","0.3333333333333333","0.03120566666666667","0.09361700000000001","0.003051451593225274","0.003051451593225274","5.975177588576889","53.776598297192","0.002975763955543237","0.002975763955543237"),

    @("MuSCs","ECs","3","1","4.138615666666666","12.415847","0.02107801645190694","0.02107801645190694","2","0.6666666666666666","0.138766","0.416298","0.01356925767068476","0.01356925767068476","0.5742991416006666","5.168692274406","0.0002860130364228578","0.0002860130364228578"),
    @("MuSCs","FAPs","3","1","4.138615666666666","12.415847","0.02107801645190694","0.02107801645190694","2","0.6666666666666666","0.05416133333333333","0.162484","0.00529617548814441","0.005296175488144411","0.2241529426608889","2.017376483948","0.0001116328740712942","0.0001116328740712942"),
    @("MuSCs","MuSCs","3","1","4.138615666666666","12.415847","0.02107801645190694","0.02107801645190694","3","1","10.002366","30.007098","0.9780831152479456","0.9780831152479456","41.395948631334","372.563537682006","0.02061605199452859","0.02061605199452859"),
    @("MuSCs","Resolving-Mac","3","1","4.138615666666666","12.415847","0.02107801645190694","0.02107801645190694","1","0.3333333333333333","0.03120566666666667","0.09361700000000001","0.003051451593225274","0.003051451593225274","0.1291482609554445","1.162334348599","6.431854688419997E-05","6.431854688419997E-05"),

    @("Resolving-Mac","ECs","1","0.3333333333333333","0.014556","0.043668","7.413387281768795E-05","7.413387281768795E-05","2","0.6666666666666666","0.138766","0.416298","0.01356925767068476","0.01356925767068476","0.002019877896","0.018178901064","1.00594162238898E-06","1.005941622388981E-06"),
    @("Resolving-Mac","FAPs","1","0.3333333333333333","0.014556","0.043668","7.413387281768795E-05","7.413387281768795E-05","2","0.6666666666666666","0.05416133333333333","0.162484","0.00529617548814441","0.005296175488144411","0.0007883723679999999","0.007095351311999999","3.926260000582541E-07","3.926260000582541E-07"),
    @("Resolving-Mac","MuSCs","1","0.3333333333333333","0.014556","0.043668","7.413387281768795E-05","7.413387281768795E-05","3","1","10.002366","30.007098","0.9780831152479456","0.9780831152479456","0.145594439496","1.310349955464","7.250908927091922E-05","7.250908927091922E-05"),
    @("Resolving-Mac","Resolving-Mac","1","0.3333333333333333","0.014556","0.043668","7.413387281768795E-05","7.413387281768795E-05","1","0.3333333333333333","0.03120566666666667","0.09361700000000001","0.003051451593225274","0.003051451593225274","0.000454229684","0.004088067156","2.262159243214937E-07","2.262159243214937E-07")
)

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$r = 2
foreach ($row in $rows) {
    $sending = $row[0]
    $target = $row[1]

    $ws.Range("A$r").Value2 = $sending
    $ws.Range("B$r").Value2 = "Ptn"
    $ws.Range("C$r").Value2 = "Ptprz1"
    $ws.Range("D$r").Value2 = $target

    for ($i = 0; $i -lt $cols.Length; $i++) {
        Set-Num "$($cols[$i])$r" $row[2 + $i]
    }

    $r = $r + 1
}
